$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.173.04'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '3.483.60'
$ws.Range("E3").Value = '  -0.51%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.31'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.59%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.66'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.387'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.54%  '
$ws.Range("D12").Value = '4.082.55'
$ws.Range("E12").Value = '  -0.49%  '
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("D15").Value = '3.486.66'
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("D16").Value = '64.176.79'
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '384.65'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.577'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("D23").Value = '3.623.27'
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.35'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000113'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.32%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.153'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.85%  '
$ws.Range("D34").Value = '3.513.32'
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.27'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("E38").Value = '  -1.15%  '
$ws.Range("E39").Value = '  -3.31%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '163.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0778'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.802'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.62'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.935'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.45%  '
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("D50").Value = '2.375.15'
$ws.Range("E51").Value = '  -2.47%  '
